$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the two runs that make up
#    "What is the longest path length from the root to a leaf node?"
#    into a single run (and drop the "_GoBack" bookmark that currently sits
#    between them).  A plain re-assignment of identical text is a no-op in
#    this engine, so we first stamp a placeholder value (forcing a real
#    edit, which also lets Word relocate "_GoBack" to wherever the next
#    real edit happens) and then set the final text.
# ---------------------------------------------------------------------------
$targetText = "What is the longest path length from the root to a leaf node?"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($paraText -eq $targetText) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne -1) {
    $p = $d.Paragraphs.Item($targetIndex)
    $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
    $rng.Text = "placeholder-text-forcing-a-real-edit"

    $p2 = $d.Paragraphs.Item($targetIndex)
    $rng2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
    $rng2.Text = $targetText
}

# ---------------------------------------------------------------------------
# 2) Replace the final (empty) paragraph of the body - the one right after
#    the last table and right before the sectPr - with the GO-ontology
#    result paragraphs, using InsertXML so the exact OOXML shape (including
#    the formatted-but-empty lead paragraph and the relocated "_GoBack"
#    bookmark) can be produced directly.
# ---------------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParagraphsXml = '<w:body ' + $wNs + '>' `
    + '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p>' `
    + '<w:p><w:r><w:t>BP Size: 29692</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' `
    + '<w:p><w:r><w:t>MF Size: 11111</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t>CC Size: 4206</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t>BP Root: [''GO:0008150'']</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t>MF Root: [''GO:0003674'']</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t>CC Root: [''GO:0005575'']</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t>BP Leaf Size: 13627</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t>MF Leaf Size: 9066</w:t></w:r></w:p>' `
    + '<w:p><w:r><w:t>CC Leaf Size: 2751</w:t></w:r></w:p>' `
    + '</w:body>'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRng = $lastPara.Range
$endRng.Collapse(0)
$null = $endRng.InsertXML($newParagraphsXml)
